# Updated sprint 3 burndown — correct a handful of Story Points / Guideline
# values on Sheet1 that were entered incorrectly, and move the active
# selection to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("B20").Value = 4
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 4
$ws.Range("B23").Value = 0

$ws.Range("C16").Select()
